$wb = $excel.ActiveWorkbook

$wsRoasts = $wb.Worksheets.Item("roasts")
$wsPlayers = $wb.Worksheets.Item("players")

$wsRoasts.Range("A5").Value = 1
$wsRoasts.Range("B5").Value = "hehehe prut"

$wsRoasts.Range("A6").Value = 1
$wsRoasts.Range("B6").Value = "du er dum og grim"

$wsRoasts.Range("A7").Value = 2
$wsRoasts.Range("B7").Value = "hej med dig"

$wsRoasts.Range("A8").Value = 2
$wsRoasts.Range("B8").Value = "sån er det bare"

$wsPlayers.Activate()
$wsPlayers.Range("F10").Select()

$wsRoasts.Activate()
$wsRoasts.Range("B8").Select()
